$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 1.45
$ws.Range("I2").Value = 6.25
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 12
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 11
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 201
$ws.Range("AJ2").Value = 19
$ws.Range("AU2").Value = 8
$ws.Range("AX2").Value = 8
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 101

# Row 3 updates
$ws.Range("H3").Value = 3.3
$ws.Range("J3").Value = 2.75
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3.2
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.73
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 8.5
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 19
$ws.Range("AL3").Value = 34
$ws.Range("AT3").Value = 2.63
$ws.Range("AV3").Value = 51

$wb.Save()
